$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source feed had swapped the fixtures now sitting in row 83 (id 81)
# and row 84 (id 82): every data column for the "Vancouver FC vs Cavalry FC"
# match and the "Forge FC vs Atletico Ottawa" match landed on the wrong
# rows. Re-sync the sheet by exchanging all data columns (B:AC) between the
# two rows while leaving the row-index column (A) untouched.
$rng83 = $ws.Range("B83:AC83")
$rng84 = $ws.Range("B84:AC84")

$vals83 = $rng83.Value2
$vals84 = $rng84.Value2

$rng83.Value2 = $vals84
$rng84.Value2 = $vals83
